$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the auto-check column (J): DEC2BIN(I<row>,8) for every data row (2-129).
# Written in chunks that mirror the existing column B fill boundary (row 66/67)
# so the resulting shared-formula groups line up the same way.
$ws.Range("J67:J129").Formula = "=DEC2BIN(I67,8)"
$ws.Range("J2").Formula = "=DEC2BIN(I2,8)"
$ws.Range("J3:J65").Formula = "=DEC2BIN(I3,8)"
$ws.Range("J66").Formula = "=DEC2BIN(I66,8)"

# Restore the view: scrolled so row 19 is at the top, with I42 as the active cell.
$av = $ws.Application.ActiveWindow
$av.ScrollRow = 19
$av.ScrollColumn = 1
$ws.Range("I42").Select()
